$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = [double]"25.8100000000006"
$ws.Cells.Item(2, 8).Value = [double]"3.157047956392489e-10"
$ws.Cells.Item(2, 9).Value = [double]"3.157047956392489e-10"
$ws.Cells.Item(2, 12).Value = [double]"53.06859653084248"
$ws.Cells.Item(2, 13).Value = "[37.47222771960861, 68.66496534207636]"
$ws.Cells.Item(2, 14).Value = [double]"1.682501582500606e-08"
$ws.Cells.Item(2, 15).Value = [double]"1.682501582500606e-08"
$ws.Cells.Item(2, 16).Value = [double]"1.805079262422964"
$ws.Cells.Item(2, 17).Value = "[1.478026573760964, 2.132131951084965]"
$ws.Cells.Item(2, 18).Value = [double]"1.709743457922741e-14"
$ws.Cells.Item(2, 19).Value = [double]"1.709743457922741e-14"
$ws.Cells.Item(2, 20).Value = [double]"60.39183246219416"
$ws.Cells.Item(2, 21).Value = "[51.234896320655054, 69.54876860373327]"
$ws.Cells.Item(2, 22).Value = [double]"0"
$ws.Cells.Item(2, 23).Value = [double]"0"
$ws.Cells.Item(2, 24).Value = [double]"18.39511511511554"
$ws.Cells.Item(2, 25).Value = [double]"17.05165165165204"
$ws.Cells.Item(2, 26).Value = [double]"19.73857857857903"
$ws.Cells.Item(3, 6).Value = [double]"25.8100000000006"
$ws.Cells.Item(3, 8).Value = [double]"3.00158168320408e-08"
$ws.Cells.Item(3, 9).Value = [double]"3.00158168320408e-08"
$ws.Cells.Item(3, 12).Value = [double]"54.43578015328404"
$ws.Cells.Item(3, 13).Value = "[35.22225158113652, 73.64930872543155]"
$ws.Cells.Item(3, 14).Value = [double]"8.502035238056038e-07"
$ws.Cells.Item(3, 15).Value = [double]"8.502035238056038e-07"
$ws.Cells.Item(3, 16).Value = [double]"1.239026532046425"
$ws.Cells.Item(3, 17).Value = "[0.8490790955648091, 1.6289739685280402]"
$ws.Cells.Item(3, 18).Value = [double]"7.950262403966235e-08"
$ws.Cells.Item(3, 19).Value = [double]"7.950262403966235e-08"
$ws.Cells.Item(3, 20).Value = [double]"62.20756606236501"
$ws.Cells.Item(3, 21).Value = "[51.30751582819373, 73.10761629653629]"
$ws.Cells.Item(3, 22).Value = [double]"5.551115123125783e-15"
$ws.Cells.Item(3, 23).Value = [double]"5.551115123125783e-15"
$ws.Cells.Item(3, 24).Value = [double]"20.72034034034082"
$ws.Cells.Item(3, 25).Value = [double]"19.11851851851896"
$ws.Cells.Item(3, 26).Value = [double]"22.32216216216268"
$ws.Cells.Item(4, 6).Value = [double]"25.8100000000006"
$ws.Cells.Item(4, 8).Value = [double]"1.235228933582633e-07"
$ws.Cells.Item(4, 9).Value = [double]"1.235228933582633e-07"
$ws.Cells.Item(4, 12).Value = [double]"60.51876588014272"
$ws.Cells.Item(4, 13).Value = "[36.23002122339051, 84.80751053689494]"
$ws.Cells.Item(4, 14).Value = [double]"8.6517036510525e-06"
$ws.Cells.Item(4, 15).Value = [double]"8.6517036510525e-06"
$ws.Cells.Item(4, 16).Value = [double]"0.8239211964369622"
$ws.Cells.Item(4, 17).Value = "[0.44655270951926873, 1.2012896833546556]"
$ws.Cells.Item(4, 18).Value = [double]"6.627649040957806e-05"
$ws.Cells.Item(4, 19).Value = [double]"6.627649040957806e-05"
$ws.Cells.Item(4, 20).Value = [double]"66.95635258428435"
$ws.Cells.Item(4, 21).Value = "[54.495799971284995, 79.4169051972837]"
$ws.Cells.Item(4, 22).Value = [double]"4.130029651605582e-14"
$ws.Cells.Item(4, 23).Value = [double]"4.130029651605582e-14"
$ws.Cells.Item(4, 24).Value = [double]"22.42550550550602"
$ws.Cells.Item(4, 25).Value = [double]"20.87535535535584"
$ws.Cells.Item(4, 26).Value = [double]"23.97565565565621"
$ws.Cells.Item(5, 6).Value = [double]"25.8100000000006"
$ws.Cells.Item(5, 8).Value = [double]"3.715696639261523e-10"
$ws.Cells.Item(5, 9).Value = [double]"3.715696639261523e-10"
$ws.Cells.Item(5, 12).Value = [double]"61.53187721288579"
$ws.Cells.Item(5, 13).Value = "[43.06792313371972, 79.99583129205186]"
$ws.Cells.Item(5, 14).Value = [double]"2.726998937774283e-08"
$ws.Cells.Item(5, 15).Value = [double]"2.726998937774283e-08"
$ws.Cells.Item(5, 16).Value = [double]"0.4213948103914236"
$ws.Cells.Item(5, 17).Value = "[0.11950002085726918, 0.723289599925578]"
$ws.Cells.Item(5, 18).Value = [double]"0.007281719018265909"
$ws.Cells.Item(5, 19).Value = [double]"0.007281719018265909"
$ws.Cells.Item(5, 20).Value = [double]"63.37362346676471"
$ws.Cells.Item(5, 21).Value = "[53.493036853194454, 73.25421008033496]"
$ws.Cells.Item(5, 22).Value = [double]"2.220446049250313e-16"
$ws.Cells.Item(5, 23).Value = [double]"2.220446049250313e-16"
$ws.Cells.Item(5, 24).Value = [double]"24.07899899899956"
$ws.Cells.Item(5, 25).Value = [double]"22.83887887887941"
$ws.Cells.Item(5, 26).Value = [double]"25.3191191191197"
$ws.Cells.Item(6, 6).Value = [double]"22.25000000000004"
$ws.Cells.Item(6, 8).Value = [double]"5.624444598950618e-08"
$ws.Cells.Item(6, 9).Value = [double]"5.624444598950618e-08"
$ws.Cells.Item(6, 12).Value = [double]"55.1315792345347"
$ws.Cells.Item(6, 13).Value = "[34.70381076475722, 75.55934770431219]"
$ws.Cells.Item(6, 14).Value = [double]"2.129497983727191e-06"
$ws.Cells.Item(6, 15).Value = [double]"2.129497983727191e-06"
$ws.Cells.Item(6, 16).Value = [double]"-0.3018947895341544"
$ws.Cells.Item(6, 17).Value = "[-0.729579074707539, 0.12578949563923025]"
$ws.Cells.Item(6, 18).Value = [double]"0.1620016016951342"
$ws.Cells.Item(6, 19).Value = [double]"0.1620016016951342"
$ws.Cells.Item(6, 20).Value = [double]"59.2782506632596"
$ws.Cells.Item(6, 21).Value = "[47.74926706212966, 70.80723426438954]"
$ws.Cells.Item(6, 22).Value = [double]"1.718625242119742e-13"
$ws.Cells.Item(6, 23).Value = [double]"1.718625242119742e-13"
$ws.Cells.Item(6, 24).Value = [double]"1.069069069069069"
$ws.Cells.Item(6, 25).Value = [double]"-0.4454454454454455"
$ws.Cells.Item(6, 26).Value = [double]"2.583583583583584"
$ws.Cells.Item(7, 6).Value = [double]"22.25000000000004"
$ws.Cells.Item(7, 8).Value = [double]"2.773692542312745e-08"
$ws.Cells.Item(7, 9).Value = [double]"2.773692542312745e-08"
$ws.Cells.Item(7, 12).Value = [double]"57.74748525500062"
$ws.Cells.Item(7, 13).Value = "[38.020434930481485, 77.47453557951975]"
$ws.Cells.Item(7, 14).Value = [double]"4.455226614297203e-07"
$ws.Cells.Item(7, 15).Value = [double]"4.455226614297203e-07"
$ws.Cells.Item(7, 16).Value = [double]"0.1572368695490383"
$ws.Cells.Item(7, 17).Value = "[-0.25786846606042424, 0.5723422051585008]"
$ws.Cells.Item(7, 18).Value = [double]"0.449493376385325"
$ws.Cells.Item(7, 19).Value = [double]"0.449493376385325"
$ws.Cells.Item(7, 20).Value = [double]"52.56679644972821"
$ws.Cells.Item(7, 21).Value = "[40.911155607788956, 64.22243729166746]"
$ws.Cells.Item(7, 22).Value = [double]"9.659606448053637e-12"
$ws.Cells.Item(7, 23).Value = [double]"9.659606448053637e-12"
$ws.Cells.Item(7, 24).Value = [double]"21.69319319319323"
$ws.Cells.Item(7, 25).Value = [double]"20.22322322322326"
$ws.Cells.Item(7, 26).Value = [double]"23.1631631631632"
$ws.Cells.Item(8, 6).Value = [double]"22.25000000000004"
$ws.Cells.Item(8, 8).Value = [double]"3.10823998761478e-08"
$ws.Cells.Item(8, 9).Value = [double]"3.10823998761478e-08"
$ws.Cells.Item(8, 12).Value = [double]"55.82618793281382"
$ws.Cells.Item(8, 13).Value = "[35.52007667504679, 76.13229919058085]"
$ws.Cells.Item(8, 14).Value = [double]"1.510266238735625e-06"
$ws.Cells.Item(8, 15).Value = [double]"1.510266238735625e-06"
$ws.Cells.Item(8, 16).Value = [double]"0.5094474573388847"
$ws.Cells.Item(8, 17).Value = "[0.10692107129334438, 0.9119738433844251]"
$ws.Cells.Item(8, 18).Value = [double]"0.01427963940752375"
$ws.Cells.Item(8, 19).Value = [double]"0.01427963940752375"
$ws.Cells.Item(8, 20).Value = [double]"56.91624963877202"
$ws.Cells.Item(8, 21).Value = "[45.80078229358483, 68.0317169839592]"
$ws.Cells.Item(8, 22).Value = [double]"1.960653861488026e-13"
$ws.Cells.Item(8, 23).Value = [double]"1.960653861488026e-13"
$ws.Cells.Item(8, 24).Value = [double]"20.44594594594598"
$ws.Cells.Item(8, 25).Value = [double]"19.02052052052055"
$ws.Cells.Item(8, 26).Value = [double]"21.87137137137142"
